$d = $word.ActiveDocument

$findText = "and the  quality measure rating."
$replaceText = "and the quality measure rating."
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 0 NOT FOUND: $findText" }
Write-Host "replacement 0 found: $found"

$findText = "For profit on the left, government in the middle of the graph and non profit on the right."
$replaceText = "On the first graph, we are simply showing the quality measure rating each nursing home received; For profit on the left, government in the middle of the graph and nonprofit on the right."
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 1 NOT FOUND: $findText" }
Write-Host "replacement 1 found: $found"

$findText = "The variance between for-profit subtypes: Corp, individual and LLC  is minimal, with partnership sub-type nearly 1 point above the others.  Similarly, there is also a minimal variance between non-profit sub types church related” and “other” while Corp is nearly 1 point higher. "
$replaceText = "The differences in the for-profit subtypes: “Corp”, “individual” and “LLC” is minimal, with “Partnership” sub-type nearly 1 point above the others.  Similarly, there is also minimal difference between nonprofit sub types “Church related” and “Other” while “Corp” is nearly 1 point higher. "
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 2 NOT FOUND: $findText" }
Write-Host "replacement 2 found: $found"

$findText = "One hypothesis could be the higher amount"
$replaceText = "One hypothesis:  Could having a higher "
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 3 NOT FOUND: $findText" }
Write-Host "replacement 3 found: $found"

$findText = "higher  of fines the lower the qm rating, so a smaller data frame was created to show the ownership type and the "
$replaceText = "higher "
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 4 NOT FOUND: $findText" }
Write-Host "replacement 4 found: $found"

$findText = " of fines each had."
$replaceText = " of fines lower the quality measure rating?  A bar graph with the only 2 government owned nursing homes in Oklahoma shows a stark contrast in fines.  The city nursing home has no fines while the county nursing home has a little over `$10,000."
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 5 NOT FOUND: $findText" }
Write-Host "replacement 5 found: $found"

$findText = "Normally, one would expect the amount of fines to decrease as the qm rating increases"
$replaceText = "Normally, one would expect the fines to decrease as the quality measure rating increases"
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
if (-not $found) { Write-Host "WARNING: replacement 6 NOT FOUND: $findText" }
Write-Host "replacement 6 found: $found"
